$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.6594644287856457
$ws.Range("C2").Value = 0.5528200562108995
$ws.Range("D2").Value = -0.5161509248614538

$ws.Range("B3").Value = -0.7420218047750701
$ws.Range("C3").Value = -0.7235027927716031
$ws.Range("D3").Value = -0.7217693428549651

$ws.Range("B4").Value = -0.6312246050722827
$ws.Range("C4").Value = 0.5532091562655042
$ws.Range("D4").Value = -0.09814806294152349

$ws.Range("B5").Value = 0.6670915658037653
$ws.Range("C5").Value = -0.6488536836452636
$ws.Range("D5").Value = -0.5740329548498322

$ws.Range("B6").Value = 0.8000244920464049
$ws.Range("C6").Value = -0.6993995071951961
$ws.Range("D6").Value = -0.6619501785003117

$ws.Range("B7").Value = -0.6538899606634785
$ws.Range("C7").Value = -0.6582497949293676
$ws.Range("D7").Value = 0.6765693485992146

$ws.Range("B8").Value = -0.8353421471694314
$ws.Range("C8").Value = 0.8846250243691259
$ws.Range("D8").Value = 0.6750595056807498

$ws.Range("B9").Value = 0.8023964163849941
$ws.Range("C9").Value = 0.7104722060248243
$ws.Range("D9").Value = 0.6738028332491574
